# Update "想去人数" (F column) figures to the refreshed scrape values
# (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value  = 1214
$ws1.Range("F5").Value  = 1163
$ws1.Range("F6").Value  = 3761
$ws1.Range("F7").Value  = 2502
$ws1.Range("F8").Value  = 62
$ws1.Range("F9").Value  = 2338
$ws1.Range("F10").Value = 243
$ws1.Range("F13").Value = 1618
$ws1.Range("F14").Value = 629
$ws1.Range("F17").Value = 21
$ws1.Range("F21").Value = 422
$ws1.Range("F22").Value = 21
$ws1.Range("F24").Value = 464
$ws1.Range("F25").Value = 655
$ws1.Range("F26").Value = 74
$ws1.Range("F28").Value = 349
$ws1.Range("F30").Value = 1603
$ws1.Range("F31").Value = 784
$ws1.Range("F32").Value = 808
$ws1.Range("F33").Value = 1904
$ws1.Range("F35").Value = 495
$ws1.Range("F37").Value = 554
$ws1.Range("F38").Value = 1186
$ws1.Range("F40").Value = 398

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F6").Value = 2
$ws2.Range("F8").Value = 6

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value  = 1214
$ws4.Range("F5").Value  = 1163
$ws4.Range("F6").Value  = 3761
$ws4.Range("F7").Value  = 2502
$ws4.Range("F8").Value  = 62
$ws4.Range("F9").Value  = 2338
$ws4.Range("F10").Value = 243
$ws4.Range("F13").Value = 1618
$ws4.Range("F14").Value = 629
$ws4.Range("F17").Value = 21
$ws4.Range("F21").Value = 422
$ws4.Range("F22").Value = 21
$ws4.Range("F24").Value = 464
$ws4.Range("F25").Value = 655
$ws4.Range("F26").Value = 74
$ws4.Range("F31").Value = 349
$ws4.Range("F33").Value = 1603
$ws4.Range("F34").Value = 784
$ws4.Range("F36").Value = 808
$ws4.Range("F37").Value = 1904
$ws4.Range("F39").Value = 2
$ws4.Range("F41").Value = 6
$ws4.Range("F42").Value = 496
$ws4.Range("F44").Value = 554
$ws4.Range("F45").Value = 1186
$ws4.Range("F47").Value = 398
